# Generate Report for Handback
# The e0d7219b-... file has now been handed back and is in sync, so:
#  - Overview sheet: the zh-cn / de-de status cells move from "Ready for
#    handoff" to "Handed back: in sync with en-US"
#  - zh-cn / de-de detail sheets: the Status cell for that row becomes
#    "Handed back: in sync with en-US", the Latest Handback DateTime is
#    refreshed, and the stale Error Detail message is cleared out.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# ---- zh-cn detail sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-09-02 12:56:38"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).ColumnWidth = 12.76

# ---- de-de detail sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-09-02 12:56:44"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).ColumnWidth = 12.76
